$d = $word.ActiveDocument

$d.Content.Find.Execute('2024-12-05 Thursday', $false, $false, $false, $false, $false, $true, 1, $false, '2024-12-06 Friday', 2) | Out-Null
$d.Content.Find.Execute('919×4=3676', $false, $false, $false, $false, $false, $true, 1, $false, '355×8=2840', 2) | Out-Null
$d.Content.Find.Execute('375×3=1125', $false, $false, $false, $false, $false, $true, 1, $false, '770×5=3850', 2) | Out-Null
$d.Content.Find.Execute('222×9=1998', $false, $false, $false, $false, $false, $true, 1, $false, '453×4=1812', 2) | Out-Null
$d.Content.Find.Execute('794×6=4764', $false, $false, $false, $false, $false, $true, 1, $false, '586×7=4102', 2) | Out-Null
$d.Content.Find.Execute('368×2=736', $false, $false, $false, $false, $false, $true, 1, $false, '922×4=3688', 2) | Out-Null
$d.Content.Find.Execute('230×3=690', $false, $false, $false, $false, $false, $true, 1, $false, '153×4=612', 2) | Out-Null
$d.Content.Find.Execute('601×6=3606', $false, $false, $false, $false, $false, $true, 1, $false, '763×8=6104', 2) | Out-Null
$d.Content.Find.Execute('641×6=3846', $false, $false, $false, $false, $false, $true, 1, $false, '272×3=816', 2) | Out-Null
$d.Content.Find.Execute('561×6=3366', $false, $false, $false, $false, $false, $true, 1, $false, '547×8=4376', 2) | Out-Null
$d.Content.Find.Execute('757×8=6056', $false, $false, $false, $false, $false, $true, 1, $false, '152×6=912', 2) | Out-Null
$d.Content.Find.Execute('463×7=3241', $false, $false, $false, $false, $false, $true, 1, $false, '881×2=1762', 2) | Out-Null
$d.Content.Find.Execute('255×5=1275', $false, $false, $false, $false, $false, $true, 1, $false, '210×8=1680', 2) | Out-Null
$d.Content.Find.Execute('603×2=1206', $false, $false, $false, $false, $false, $true, 1, $false, '774×3=2322', 2) | Out-Null
$d.Content.Find.Execute('797×6=4782', $false, $false, $false, $false, $false, $true, 1, $false, '190×5=950', 2) | Out-Null
$d.Content.Find.Execute('689×6=4134', $false, $false, $false, $false, $false, $true, 1, $false, '614×8=4912', 2) | Out-Null
$d.Content.Find.Execute('415×6=2490', $false, $false, $false, $false, $false, $true, 1, $false, '913×9=8217', 2) | Out-Null
$d.Content.Find.Execute('257×3=771', $false, $false, $false, $false, $false, $true, 1, $false, '232×8=1856', 2) | Out-Null
$d.Content.Find.Execute('279×2=558', $false, $false, $false, $false, $false, $true, 1, $false, '663×9=5967', 2) | Out-Null
$d.Content.Find.Execute('889×2=1778', $false, $false, $false, $false, $false, $true, 1, $false, '722×2=1444', 2) | Out-Null
$d.Content.Find.Execute('137×4=548', $false, $false, $false, $false, $false, $true, 1, $false, '604×4=2416', 2) | Out-Null
$d.Content.Find.Execute('943×8=7544', $false, $false, $false, $false, $false, $true, 1, $false, '649×8=5192', 2) | Out-Null
$d.Content.Find.Execute('851×8=6808', $false, $false, $false, $false, $false, $true, 1, $false, '899×7=6293', 2) | Out-Null
$d.Content.Find.Execute('568×3=1704', $false, $false, $false, $false, $false, $true, 1, $false, '902×5=4510', 2) | Out-Null
$d.Content.Find.Execute('886×5=4430', $false, $false, $false, $false, $false, $true, 1, $false, '320×4=1280', 2) | Out-Null
$d.Content.Find.Execute('914×5=4570', $false, $false, $false, $false, $false, $true, 1, $false, '450×4=1800', 2) | Out-Null
